# Add 2022-Q1 data
# 1) Insert a new worksheet "2022-Q1" right before the "总计" (total) sheet,
#    copying the per-fund-holding layout/format used by the other quarter sheets.
# 2) Fill the new sheet with the 2022-Q1 fund holdings data.
# 3) Prepend a 2022-Q1 summary row to the "总计" sheet and renumber the index column.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item("总计")

$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# NOTE: the sheet reference passed as the "insert before" target above
# becomes stale once the insertion happens, so re-fetch "总计" by name
# before using it for any further reads/writes.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy the header + first-row formatting from the template sheet, then stamp
# the same per-row formatting (border/alignment on column A) down to row 10.
$template.Range("A1:H2").Copy()
$newSheet.Range("A1:H2").PasteSpecial(-4122)
$template.Range("A2:H2").Copy()
$newSheet.Range("A3:H10").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 - 003318
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'003318"
$newSheet.Range("C2").Value = "景顺长城中证500行业中性低波动指数"
$newSheet.Range("D2").Value = "'13.99"
$newSheet.Range("E2").Value = "'93.88"
$newSheet.Range("F2").Value = "'1.91"
$newSheet.Range("G2").Value = "'0.2672"
$newSheet.Range("H2").Value = 2

# Row 3 - 512330
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'512330"
$newSheet.Range("C3").Value = "南方中证500信息技术指数ETF"
$newSheet.Range("D3").Value = "'5.09"
$newSheet.Range("E3").Value = "'99.74"
$newSheet.Range("F3").Value = "'2.81"
$newSheet.Range("G3").Value = "'0.1430"
$newSheet.Range("H3").Value = 9

# Row 4 - 005994
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'005994"
$newSheet.Range("C4").Value = "国投瑞银中证500指数量化增强A"
$newSheet.Range("D4").Value = "'11.53"
$newSheet.Range("E4").Value = "'87.00"
$newSheet.Range("F4").Value = "'1.22"
$newSheet.Range("G4").Value = "'0.1407"
$newSheet.Range("H4").Value = 7

# Row 5 - 002311
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'002311"
$newSheet.Range("C5").Value = "创金合信中证500指数增强A"
$newSheet.Range("D5").Value = "'5.72"
$newSheet.Range("E5").Value = "'92.95"
$newSheet.Range("F5").Value = "'1.09"
$newSheet.Range("G5").Value = "'0.0623"
$newSheet.Range("H5").Value = 8

# Row 6 - 007089
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'007089"
$newSheet.Range("C6").Value = "国投瑞银中证500指数量化增强C"
$newSheet.Range("D6").Value = "'3.82"
$newSheet.Range("E6").Value = "'87.00"
$newSheet.Range("F6").Value = "'1.22"
$newSheet.Range("G6").Value = "'0.0466"
$newSheet.Range("H6").Value = 7

# Row 7 - 002316
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'002316"
$newSheet.Range("C7").Value = "创金合信中证500指数增强C"
$newSheet.Range("D7").Value = "'2.60"
$newSheet.Range("E7").Value = "'92.95"
$newSheet.Range("F7").Value = "'1.09"
$newSheet.Range("G7").Value = "'0.0283"
$newSheet.Range("H7").Value = 8

# Row 8 - 512260
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "'512260"
$newSheet.Range("C8").Value = "华安中证500行业中性低波动ETF"
$newSheet.Range("D8").Value = "'1.17"
$newSheet.Range("E8").Value = "'96.94"
$newSheet.Range("F8").Value = "'1.98"
$newSheet.Range("G8").Value = "'0.0232"
$newSheet.Range("H8").Value = 2

# Row 9 - 004192
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "'004192"
$newSheet.Range("C9").Value = "招商中证500指数增强A"
$newSheet.Range("D9").Value = "'0.96"
$newSheet.Range("E9").Value = "'94.32"
$newSheet.Range("F9").Value = "'1.00"
$newSheet.Range("G9").Value = "'0.0096"
$newSheet.Range("H9").Value = 7

# Row 10 - 004193
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "'004193"
$newSheet.Range("C10").Value = "招商中证500指数增强C"
$newSheet.Range("D10").Value = "'0.42"
$newSheet.Range("E10").Value = "'94.32"
$newSheet.Range("F10").Value = "'1.00"
$newSheet.Range("G10").Value = "'0.0042"
$newSheet.Range("H10").Value = 7

# Update the "总计" (grand total) sheet: add a 2022-Q1 row at the top of the
# data area and renumber/rewrite the remaining rows underneath it.
$totalSheet.Range("A2:D10").ClearContents()
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 0.73

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 5
$totalSheet.Range("D3").Value = 0.71

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 5
$totalSheet.Range("D4").Value = 0.37

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 31
$totalSheet.Range("D5").Value = 6

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q1"
$totalSheet.Range("C6").Value = 8
$totalSheet.Range("D6").Value = 12.69

$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2020-Q4"
$totalSheet.Range("C7").Value = 10
$totalSheet.Range("D7").Value = 5.71
